$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Carslaw building"
$ws.Range("B6").Value = "My university office"

$ws.Range("B6").Select()
